$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.358.17"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.718.13"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.93"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4725"
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2637"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06208"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.717.31"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07062"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.31"
$ws.Range("E12").Value = "  +2.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5916"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.409"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.16"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.354.68"
$ws.Range("E18").Value = "  +3.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006802"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.56"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.937.39"
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.547"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.756"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.322"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.23"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.25"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.410"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "108.06"
$ws.Range("E28").Value = "  +2.94%  "
$ws.Range("E29").Value = "  +3.74%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.689"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07742"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04434"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9792"
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6187"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9300"
$ws.Range("E37").Value = "  +7.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "113.35"
$ws.Range("E38").Value = "  +16.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.412"
$ws.Range("E39").Value = "  -7.69%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01478"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.340"
$ws.Range("E43").Value = "  +13.40%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1168"
$ws.Range("E45").Value = "  +4.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.297"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.43"
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.686"
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3371"
$ws.Range("E51").Value = "  +0.78%  "
